$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "-"
$ws.Range("D3").Value = "[-, -, -, 'MCT-3A-Tecnologia da soldagem']"
$ws.Range("E3").Value = "-"
$ws.Range("D4").Value = "[-, -, -, 'MCT-3A-Tecnologia da soldagem']"
$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = "[-, -, -, 'MCT-3A-Tecnologia da soldagem']"
$ws.Range("B7").Value = "-"
$ws.Range("D8").Value = "[-, -, -, 'MCT-3A-Tecnologia da soldagem']"
